$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt2"
$ws.Cells.Item(2, 3).Value = "Fzd5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01070233333333333
$ws.Cells.Item(2, 8).Value = 0.032107
$ws.Cells.Item(2, 9).Value = 0.006017198313602724
$ws.Cells.Item(2, 10).Value = 0.006017198313602724
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.343367666666667
$ws.Cells.Item(2, 14).Value = 7.030103
$ws.Cells.Item(2, 15).Value = 0.07041630712807322
$ws.Cells.Item(2, 16).Value = 0.07041630712807322
$ws.Cells.Item(2, 17).Value = 0.02507950189122222
$ws.Cells.Item(2, 18).Value = 0.225715517021
$ws.Cells.Item(2, 19).Value = 0.0004237088845011737
$ws.Cells.Item(2, 20).Value = 0.0004237088845011737

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt2"
$ws.Cells.Item(3, 3).Value = "Fzd5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01070233333333333
$ws.Cells.Item(3, 8).Value = 0.032107
$ws.Cells.Item(3, 9).Value = 0.006017198313602724
$ws.Cells.Item(3, 10).Value = 0.006017198313602724
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.842594333333333
$ws.Cells.Item(3, 14).Value = 14.527783
$ws.Cells.Item(3, 15).Value = 0.1455160514174545
$ws.Cells.Item(3, 16).Value = 0.1455160514174545
$ws.Cells.Item(3, 17).Value = 0.05182705875344444
$ws.Cells.Item(3, 18).Value = 0.4664435287809999
$ws.Cells.Item(3, 19).Value = 0.0008755989391912343
$ws.Cells.Item(3, 20).Value = 0.0008755989391912343

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt2"
$ws.Cells.Item(4, 3).Value = "Fzd5"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.01070233333333333
$ws.Cells.Item(4, 8).Value = 0.032107
$ws.Cells.Item(4, 9).Value = 0.006017198313602724
$ws.Cells.Item(4, 10).Value = 0.006017198313602724
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.196962666666667
$ws.Cells.Item(4, 14).Value = 12.590888
$ws.Cells.Item(4, 15).Value = 0.1261153409022843
$ws.Cells.Item(4, 16).Value = 0.1261153409022843
$ws.Cells.Item(4, 17).Value = 0.04491729344622222
$ws.Cells.Item(4, 18).Value = 0.4042556410159999
$ws.Cells.Item(4, 19).Value = 0.0007588610165966577
$ws.Cells.Item(4, 20).Value = 0.0007588610165966577

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt2"
$ws.Cells.Item(5, 3).Value = "Fzd5"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.01070233333333333
$ws.Cells.Item(5, 8).Value = 0.032107
$ws.Cells.Item(5, 9).Value = 0.006017198313602724
$ws.Cells.Item(5, 10).Value = 0.006017198313602724
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.056184333333333
$ws.Cells.Item(5, 14).Value = 15.168553
$ws.Cells.Item(5, 15).Value = 0.151934258535964
$ws.Cells.Item(5, 16).Value = 0.151934258535964
$ws.Cells.Item(5, 17).Value = 0.05411297013011111
$ws.Cells.Item(5, 18).Value = 0.4870167311709999
$ws.Cells.Item(5, 19).Value = 0.0009142185642410832
$ws.Cells.Item(5, 20).Value = 0.0009142185642410832

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Wnt2"
$ws.Cells.Item(6, 3).Value = "Fzd5"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.01070233333333333
$ws.Cells.Item(6, 8).Value = 0.032107
$ws.Cells.Item(6, 9).Value = 0.006017198313602724
$ws.Cells.Item(6, 10).Value = 0.006017198313602724
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 15.01164366666667
$ws.Cells.Item(6, 14).Value = 45.034931
$ws.Cells.Item(6, 15).Value = 0.451087776777607
$ws.Cells.Item(6, 16).Value = 0.4510877767776071
$ws.Cells.Item(6, 17).Value = 0.1606596144018889
$ws.Cells.Item(6, 18).Value = 1.445936529617
$ws.Cells.Item(6, 19).Value = 0.002714284609713019
$ws.Cells.Item(6, 20).Value = 0.00271428460971302

# Row 7
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Wnt2"
$ws.Cells.Item(7, 3).Value = "Fzd5"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.01070233333333333
$ws.Cells.Item(7, 8).Value = 0.032107
$ws.Cells.Item(7, 9).Value = 0.006017198313602724
$ws.Cells.Item(7, 10).Value = 0.006017198313602724
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.828011333333333
$ws.Cells.Item(7, 14).Value = 5.484034
$ws.Cells.Item(7, 15).Value = 0.05493026523861683
$ws.Cells.Item(7, 16).Value = 0.05493026523861684
$ws.Cells.Item(7, 17).Value = 0.01956398662644444
$ws.Cells.Item(7, 18).Value = 0.176075879638
$ws.Cells.Item(7, 19).Value = 0.0003305262993595556
$ws.Cells.Item(7, 20).Value = 0.0003305262993595556

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt2"
$ws.Cells.Item(8, 3).Value = "Fzd5"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.767921666666667
$ws.Cells.Item(8, 8).Value = 5.303765
$ws.Cells.Item(8, 9).Value = 0.9939828016863973
$ws.Cells.Item(8, 10).Value = 0.9939828016863973
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.343367666666667
$ws.Cells.Item(8, 14).Value = 7.030103
$ws.Cells.Item(8, 15).Value = 0.07041630712807322
$ws.Cells.Item(8, 16).Value = 0.07041630712807322
$ws.Cells.Item(8, 17).Value = 4.142890470866111
$ws.Cells.Item(8, 18).Value = 37.286014237795
$ws.Cells.Item(8, 19).Value = 0.06999259824357204
$ws.Cells.Item(8, 20).Value = 0.06999259824357204

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt2"
$ws.Cells.Item(9, 3).Value = "Fzd5"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.767921666666667
$ws.Cells.Item(9, 8).Value = 5.303765
$ws.Cells.Item(9, 9).Value = 0.9939828016863973
$ws.Cells.Item(9, 10).Value = 0.9939828016863973
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 4.842594333333333
$ws.Cells.Item(9, 14).Value = 14.527783
$ws.Cells.Item(9, 15).Value = 0.1455160514174545
$ws.Cells.Item(9, 16).Value = 0.1455160514174545
$ws.Cells.Item(9, 17).Value = 8.561327444777223
$ws.Cells.Item(9, 18).Value = 77.051947002995
$ws.Cells.Item(9, 19).Value = 0.1446404524782632
$ws.Cells.Item(9, 20).Value = 0.1446404524782632

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Wnt2"
$ws.Cells.Item(10, 3).Value = "Fzd5"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.767921666666667
$ws.Cells.Item(10, 8).Value = 5.303765
$ws.Cells.Item(10, 9).Value = 0.9939828016863973
$ws.Cells.Item(10, 10).Value = 0.9939828016863973
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.196962666666667
$ws.Cells.Item(10, 14).Value = 12.590888
$ws.Cells.Item(10, 15).Value = 0.1261153409022843
$ws.Cells.Item(10, 16).Value = 0.1261153409022843
$ws.Cells.Item(10, 17).Value = 7.419901232591112
$ws.Cells.Item(10, 18).Value = 66.77911109332
$ws.Cells.Item(10, 19).Value = 0.1253564798856876
$ws.Cells.Item(10, 20).Value = 0.1253564798856876

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Wnt2"
$ws.Cells.Item(11, 3).Value = "Fzd5"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.767921666666667
$ws.Cells.Item(11, 8).Value = 5.303765
$ws.Cells.Item(11, 9).Value = 0.9939828016863973
$ws.Cells.Item(11, 10).Value = 0.9939828016863973
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 5.056184333333333
$ws.Cells.Item(11, 14).Value = 15.168553
$ws.Cells.Item(11, 15).Value = 0.151934258535964
$ws.Cells.Item(11, 16).Value = 0.151934258535964
$ws.Cells.Item(11, 17).Value = 8.938937833560557
$ws.Cells.Item(11, 18).Value = 80.450440502045
$ws.Cells.Item(11, 19).Value = 0.1510200399717229
$ws.Cells.Item(11, 20).Value = 0.1510200399717229

# Row 12
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Wnt2"
$ws.Cells.Item(12, 3).Value = "Fzd5"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.767921666666667
$ws.Cells.Item(12, 8).Value = 5.303765
$ws.Cells.Item(12, 9).Value = 0.9939828016863973
$ws.Cells.Item(12, 10).Value = 0.9939828016863973
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 15.01164366666667
$ws.Cells.Item(12, 14).Value = 45.034931
$ws.Cells.Item(12, 15).Value = 0.451087776777607
$ws.Cells.Item(12, 16).Value = 0.4510877767776071
$ws.Cells.Item(12, 17).Value = 26.53941009057944
$ws.Cells.Item(12, 18).Value = 238.854690815215
$ws.Cells.Item(12, 19).Value = 0.448373492167894
$ws.Cells.Item(12, 20).Value = 0.4483734921678941

# Row 13
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Wnt2"
$ws.Cells.Item(13, 3).Value = "Fzd5"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.767921666666667
$ws.Cells.Item(13, 8).Value = 5.303765
$ws.Cells.Item(13, 9).Value = 0.9939828016863973
$ws.Cells.Item(13, 10).Value = 0.9939828016863973
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 1.828011333333333
$ws.Cells.Item(13, 14).Value = 5.484034
$ws.Cells.Item(13, 15).Value = 0.05493026523861683
$ws.Cells.Item(13, 16).Value = 0.05493026523861684
$ws.Cells.Item(13, 17).Value = 3.231780843112223
$ws.Cells.Item(13, 18).Value = 29.08602758801
$ws.Cells.Item(13, 19).Value = 0.05459973893925728
$ws.Cells.Item(13, 20).Value = 0.05459973893925728
